$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 91.34108066666666
$ws.Range("H2").Value = 274.023242
$ws.Range("I2").Value = 0.2190334467302001
$ws.Range("J2").Value = 0.2190334467302
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1575256666666667
$ws.Range("N2").Value = 0.472577
$ws.Range("O2").Value = 0.6985926944284299
$ws.Range("P2").Value = 0.69859269442843
$ws.Range("Q2").Value = 14.38856462607044
$ws.Range("R2").Value = 129.497081634634
$ws.Range("S2").Value = 0.1530151657211964
$ws.Range("T2").Value = 0.1530151657211964

$ws.Range("G3").Value = 91.34108066666666
$ws.Range("H3").Value = 274.023242
$ws.Range("I3").Value = 0.2190334467302001
$ws.Range("J3").Value = 0.2190334467302
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.06796433333333333
$ws.Range("N3").Value = 0.203893
$ws.Range("O3").Value = 0.30140730557157
$ws.Range("P3").Value = 0.30140730557157
$ws.Range("Q3").Value = 6.207935653456222
$ws.Range("R3").Value = 55.87142088110599
$ws.Range("S3").Value = 0.06601828100900362
$ws.Range("T3").Value = 0.0660182810090036

$ws.Range("G4").Value = 276.4348856666666
$ws.Range("H4").Value = 829.3046569999999
$ws.Range("I4").Value = 0.6628833966285105
$ws.Range("J4").Value = 0.6628833966285105
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1575256666666667
$ws.Range("N4").Value = 0.472577
$ws.Range("O4").Value = 0.6985926944284299
$ws.Range("P4").Value = 0.69859269442843
$ws.Range("Q4").Value = 43.54558965456544
$ws.Range("R4").Value = 391.910306891089
$ws.Range("S4").Value = 0.4630854981425807
$ws.Range("T4").Value = 0.4630854981425808

$ws.Range("G5").Value = 276.4348856666666
$ws.Range("H5").Value = 829.3046569999999
$ws.Range("I5").Value = 0.6628833966285105
$ws.Range("J5").Value = 0.6628833966285105
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.06796433333333333
$ws.Range("N5").Value = 0.203893
$ws.Range("O5").Value = 0.30140730557157
$ws.Range("P5").Value = 0.30140730557157
$ws.Range("Q5").Value = 18.78771271441122
$ws.Range("R5").Value = 169.089414429701
$ws.Range("S5").Value = 0.1997978984859297
$ws.Range("T5").Value = 0.1997978984859297

$ws.Range("G6").Value = 49.24290466666667
$ws.Range("H6").Value = 147.728714
$ws.Range("I6").Value = 0.1180831566412894
$ws.Range("J6").Value = 0.1180831566412894
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1575256666666667
$ws.Range("N6").Value = 0.472577
$ws.Range("O6").Value = 0.6985926944284299
$ws.Range("P6").Value = 0.69859269442843
$ws.Range("Q6").Value = 7.757021386219779
$ws.Range("R6").Value = 69.813192475978
$ws.Range("S6").Value = 0.08249203056465274
$ws.Range("T6").Value = 0.08249203056465274

$ws.Range("G7").Value = 49.24290466666667
$ws.Range("H7").Value = 147.728714
$ws.Range("I7").Value = 0.1180831566412894
$ws.Range("J7").Value = 0.1180831566412894
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.06796433333333333
$ws.Range("N7").Value = 0.203893
$ws.Range("O7").Value = 0.30140730557157
$ws.Range("P7").Value = 0.30140730557157
$ws.Range("Q7").Value = 3.346761187066889
$ws.Range("R7").Value = 30.120850683602
$ws.Range("S7").Value = 0.0355911260766367
$ws.Range("T7").Value = 0.03559112607663669

